# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column D) for the file rows whose
# handoff file name is shared with row 2/3 (i.e. rows 4, 6, 7 and 8 of the
# per-locale status sheets) to reflect the latest handoff run timestamp.

$wb = $excel.ActiveWorkbook

$updates = @{
    "zh-cn" = "2016-02-18 02:32:56"
    "de-de" = "2016-02-18 02:33:07"
}

$rows = @(4, 6, 7, 8)

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newValue = $updates[$sheetName]
    foreach ($r in $rows) {
        $ws.Cells.Item($r, 4).Value = $newValue
    }
}
